$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row styling for rows 7, 293, 294 (shift from banded/grey style s=3 to s=2) ---
# Copy full-row formatting from row 2 (style index 2) onto rows 7, 293, 294
$ws.Range("A2:Q2").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$ws.Range("A2:Q2").Copy()
$ws.Range("A293:Q293").PasteSpecial(-4122)
$ws.Range("A2:Q2").Copy()
$ws.Range("A294:Q294").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Reorder comma-separated IA Control identifiers in column A ---
$ws.Range("A2").Value = 'AC-7 a,AC-7 b'
$ws.Range("A3").Value = 'AC-7 a,AC-7 b'
$ws.Range("A4").Value = 'AC-7 a,AC-7 b'
$ws.Range("A5").Value = 'AC-7 a,AC-7 b'
$ws.Range("A13").Value = 'CM-6 b,CM-5 (1),AU-6 (4),AU-7 (1),AU-7 a,AU-3 (1),AU-12 a,MA-4 (1) (a),AU-14 (1),AU-3'
$ws.Range("A14").Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Range("A15").Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Range("A17").Value = 'CM-6 b,CM-7 (2)'
$ws.Range("A22").Value = 'CM-6 b,CM-7 (2)'
$ws.Range("A23").Value = 'CM-6 b,CM-7 (2)'
$ws.Range("A37").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-14 (1),AU-3,AU-12 c'
$ws.Range("A38").Value = 'AU-14 (1),AU-4'
$ws.Range("A39").Value = 'AU-4 (1),AU-3'
$ws.Range("A44").Value = 'AU-6 (4),CM-6 b,AU-4 (1)'
$ws.Range("A45").Value = 'CM-6 b,AU-4 (1)'
$ws.Range("A46").Value = 'AU-8 b,AU-8 (1) (b),AU-8 (1) (a)'
$ws.Range("A48").Value = 'IA-2 (12),IA-2 (11)'
$ws.Range("A49").Value = 'IA-2 (12),IA-2 (1),IA-2 (11)'
$ws.Range("A50").Value = 'CM-3 (5),SI-6 d,SI-6 b'
$ws.Range("A51").Value = 'CM-3 (5),SI-6 d'
$ws.Range("A52").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A53").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A54").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A55").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A56").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A57").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A58").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A59").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A60").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A61").Value = 'CM-6 b,CM-5 (1),AU-12 (3),AU-7 a,AU-12 a,AU-8 b,AU-7 b,AU-12 c'
$ws.Range("A62").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A63").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A64").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A65").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A66").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A67").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A68").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A70").Value = 'SI-11 b,AU-9'
$ws.Range("A71").Value = 'SI-11 b,AU-9'
$ws.Range("A72").Value = 'SI-11 b,AU-9'
$ws.Range("A73").Value = 'SI-11 b,AU-9'
$ws.Range("A74").Value = 'SI-11 b,AU-9'
$ws.Range("A82").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A83").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A84").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A85").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A86").Value = 'SC-13,MA-4 (6)'
$ws.Range("A87").Value = 'MA-4 (6),AC-17 (2)'
$ws.Range("A88").Value = 'SC-13,MA-4 (6)'
$ws.Range("A92").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A93").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A94").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A95").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A96").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A97").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A98").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A99").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A100").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A101").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A102").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A103").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A113").Value = 'SC-13,SC-8,MA-4 c,AC-17 (2)'
$ws.Range("A115").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A117").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A122").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A123").Value = 'SC-28,SC-28 (1)'
$ws.Range("A133").Value = 'CM-6 b,AC-6 (10)'
$ws.Range("A134").Value = 'CM-6 b,AC-6 (10)'
$ws.Range("A135").Value = 'CM-6 b,AC-6 (10)'
$ws.Range("A140").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A141").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A142").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A152").Value = 'CM-6 b,IA-7'
$ws.Range("A153").Value = 'CM-6 b,IA-7'
$ws.Range("A154").Value = 'CM-6 b,IA-7'
$ws.Range("A156").Value = 'CM-7 a,IA-7'
$ws.Range("A157").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A160").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A161").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A162").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A163").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A164").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A171").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A175").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A176").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A177").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A178").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A179").Value = 'AU-3 (1),AU-3,MA-4 (1) (a)'
$ws.Range("A180").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A181").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A182").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A183").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A184").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A185").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a)'
$ws.Range("A186").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A187").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A188").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A189").Value = 'AU-3 (1),AU-3,AU-12 c,MA-4 (1) (a)'
$ws.Range("A190").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A191").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A192").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A193").Value = 'AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A194").Value = 'AC-2 (4),AU-3 (1),AU-12 a,MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A195").Value = 'AC-2 (4),AU-3 (1),MA-4 (1) (a),AU-3,AU-12 c'
$ws.Range("A196").Value = 'IA-2 (3),IA-2 (4),IA-2 (1),IA-2 (2)'
$ws.Range("A197").Value = 'IA-2 (3),IA-2 (5),IA-2 (4),IA-2 (2),IA-2'
$ws.Range("A198").Value = 'IA-2 (3),IA-2 (5),IA-2 (4),IA-2 (2),IA-2'
$ws.Range("A201").Value = 'IA-11,AC-3 (4)'
$ws.Range("A206").Value = 'SC-8 (1),SC-8,SC-8 (2)'
$ws.Range("A207").Value = 'SC-8 (1),SC-8,SC-8 (2)'
$ws.Range("A208").Value = 'SC-8 (1),SC-8,AC-18 (1)'
$ws.Range("A213").Value = 'CM-7 a,IA-5 (1) (c),CM-6 b'
$ws.Range("A215").Value = 'CM-6 b,AU-12 a'
$ws.Range("A216").Value = 'CM-6 b,SC-5 (2),SC-5'
$ws.Range("A219").Value = 'CM-6 b,SI-16'
$ws.Range("A220").Value = 'AU-3 (1),IA-2,IA-8'
$ws.Range("A232").Value = 'SC-2,CM-6 b,SI-16'
$ws.Range("A236").Value = 'CM-7 a,AC-18 (1)'
$ws.Range("A243").Value = 'CM-6 b,IA-5 (1) (a),IA-5 (1) (b)'
$ws.Range("A250").Value = 'CM-6 b,SC-4'
$ws.Range("A253").Value = 'IA-2 (1),IA-2 (11)'
$ws.Range("A254").Value = 'IA-2 (12),IA-2 (11)'
$ws.Range("A258").Value = 'SI-6 a,SC-3'
$ws.Range("A267").Value = 'CM-6 b,CM-5 (3)'
$ws.Range("A275").Value = 'CM-6 b,CM-7 a'
$ws.Range("A339").Value = 'CM-6 b,IA-5 (1) (c)'
$ws.Range("A341").Value = 'CM-6 b,IA-2 (2)'
$ws.Range("A342").Value = 'CM-6 b,CM-5 (1)'
$ws.Range("A343").Value = 'CM-6 b,CM-5 (1)'
$ws.Range("A351").Value = 'CM-6 b,AC-17 (2)'
$ws.Range("A374").Value = 'CM-6 b,AU-3'
$ws.Range("A377").Value = 'CM-6 b,SC-3'
$ws.Range("A382").Value = 'AC-17 (1),AC-17 (9),CM-7 b,CM-6 b'
$ws.Range("A383").Value = 'AC-17 (1),CM-7 b,CM-6 b'
$ws.Range("A412").Value = 'CM-6 b,IA-3'
$ws.Range("A413").Value = 'CM-6 b,IA-3'
$ws.Range("A422").Value = 'CM-6 b,SC-3'
$ws.Range("A429").Value = 'CM-6 b,SC-3'
$ws.Range("A435").Value = 'CM-6 b,AU-4'
$ws.Range("A439").Value = 'CM-6 b,IA-3'
$ws.Range("A440").Value = 'CM-6 b,IA-3'
$ws.Range("A454").Value = 'CM-6 b,SI-2 (2)'
$ws.Range("A458").Value = 'MA-4 e,SC-10,MA-4 (7),AC-12'
$ws.Range("A462").Value = 'SC-8 (1),SC-8,AC-17 (2)'
$ws.Range("A479").Value = 'CM-7 b,IA-3'
$ws.Range("A497").Value = 'AU-4 (1),AU-4'
$ws.Range("A500").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A501").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A502").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A503").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A518").Value = 'SC-8,SC-8 (2)'
$ws.Range("A533").Value = 'SC-8,AC-17 (2)'

# --- Special-case rows with additional content changes ---
$ws.Range("A7").Value = 'AU-5 a,AU-5 b'
$ws.Range("A294").Value = 'AU-5 a,AU-5 (1)'

# --- Row 7: new Fix (M) text ---
$m7 = @'
Configure Red Hat Enterprise Linux 9 to shutdown when auditing failures occur.
If the auditd daemon is configured to use the augenrules program to read
audit rules during daemon startup (the default), add the following line to
the bottom of "/etc/audit/rules.d/immutable.rules":
-f 2
If the auditd daemon is configured to use the auditctl utility to read
audit rules during daemon startup, add the following line to the
bottom of the /etc/audit/audit.rules file:
-f 2
'@
$ws.Range("M7").Value = $m7

# --- Row 293: updated Check (K) text and new Fix (M) text ---
$k293 = @'
Find the list of alias maps used by the Postfix mail server:
 $ sudo postconf alias_maps 
Query the Postfix alias maps for an alias for the  root  user:
 $ sudo postmap -q root hash:/etc/aliases 
The output should return an alias.

If the alias is not set, then this is a finding.
'@
$ws.Range("K293").Value = $k293

$m293 = @'
Configure a valid email address as an alias for the root account.
Append the following line to "/etc/aliases":
root: system.administrator@mail.mil
Then, run the following command:
$ sudo newaliases
'@
$ws.Range("M293").Value = $m293

# --- Row 294: new Fix (M) text ---
$m294 = @'
Configure "auditd" service to notify the SA and ISSO in the event of an audit processing failure.
Edit the following line in "/etc/audit/auditd.conf" to ensure that administrators are notified via email for those situations:
action_mail_acct = root
'@
$ws.Range("M294").Value = $m294

Write-Host "Edit complete"
